# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B (TB), C (d2S), D (K), E (IP) for rows 2-8.
# Column G (sum) is recomputed as B+C+D+E for each row.
$data = @{
    2 = @{ B = 1.459612070389937;  C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732 }
    3 = @{ B = 1.459612070389937;  C = 0.3127903958511391; D = 0.1575252929769615; E = 0.496779210170732 }
    4 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732 }
    5 = @{ B = 0.3048080303191223; C = 0.3127903958511391; D = 0.1575252929769615; E = 8.660232485948974 }
    6 = @{ B = 0.6753301551942219; C = 1.667794583268128;  D = 0.1575252929769615; E = 0.496779210170732 }
    7 = @{ B = 0.01514828764759746;C = 0.3127903958511391; D = 0.1575252929769615; E = 0.496779210170732 }
    8 = @{ B = 1.459612070389937;  C = 0.3127903958511391; D = 0.8054896365839992; E = 8.660232485948974 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.B + $vals.C + $vals.D + $vals.E
}
